# Apply crypto price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (values keep their original text formatting).
$ws.Range("D4:D5").NumberFormat = "@"
$ws.Range("D7:D11").NumberFormat = "@"
$ws.Range("D13:D16").NumberFormat = "@"
$ws.Range("D18:D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24:D25").NumberFormat = "@"
$ws.Range("D27:D31").NumberFormat = "@"
$ws.Range("D33:D36").NumberFormat = "@"
$ws.Range("D38:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.971.75"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "1.958.21"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "243.97"
$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "0.4844"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").Value = "0.2941"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.07099"
$ws.Range("E9").Value = "  +4.58%  "

$ws.Range("D10").Value = "19.69"
$ws.Range("E10").Value = "  +3.21%  "

$ws.Range("D11").Value = "106.97"
$ws.Range("E11").Value = "  +1.13%  "

$ws.Range("D12").Value = "1.955.16"
$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("D13").Value = "0.07761"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").Value = "5.389"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").Value = "0.7051"
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").Value = "278.59"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Value = "30.985.58"
$ws.Range("E17").Value = "  -0.14%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000007834"
$ws.Range("E18").Value = "  +1.64%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "13.33"
$ws.Range("E19").Value = "  +1.12%  "

$ws.Range("D20").Value = "2.248.88"
$ws.Range("E20").Value = "  +1.35%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "5.519"
$ws.Range("E22").Value = "  -1.24%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "6.517"
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("D25").Value = "9.768"
$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").Value = "19.72"
$ws.Range("E27").Value = "  -1.28%  "

$ws.Range("D28").Value = "2.179"
$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").Value = "0.1051"
$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("D30").Value = "1.408"
$ws.Range("E30").Value = "  -2.81%  "

$ws.Range("D31").Value = "4.622"
$ws.Range("E31").Value = "  -3.56%  "

$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("D33").Value = "4.413"
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").Value = "0.04896"
$ws.Range("E34").Value = "  -3.52%  "

$ws.Range("D35").Value = "0.7537"
$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("D36").Value = "1.171"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").Value = "0.02008"
$ws.Range("E38").Value = "  -1.00%  "

$ws.Range("D39").Value = "2.681"
$ws.Range("E39").Value = "  -1.45%  "

$ws.Range("D40").Value = "78.35"
$ws.Range("E40").Value = "  +9.27%  "

$ws.Range("D41").Value = "6.515"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").Value = "2.127"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("D43").Value = "0.8958"
$ws.Range("E43").Value = "  +0.69%  "

$ws.Range("D44").Value = "109.44"
$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("D45").Value = "0.4451"
$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("D46").Value = "7.905"
$ws.Range("E46").Value = "  +5.62%  "

$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").Value = "984.74"
$ws.Range("E48").Value = "  +2.66%  "

$ws.Range("D49").Value = "0.1249"
$ws.Range("E49").Value = "  -1.52%  "

$ws.Range("D50").Value = "9.314"
$ws.Range("E50").Value = "  -1.30%  "

$ws.Range("D51").Value = "35.98"
$ws.Range("E51").Value = "  +0.29%  "

